$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new PI-form row (row 8) for Kevin Kiwi, only billed for
# his first entire quarter.
$ws.Range("A8").Value = 44089.3333333333
$ws.Range("B8").Value = "kkiwi@example.com"
$ws.Range("C8").Value = "Kevin"
$ws.Range("D8").Value = "Kiwi"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Yes"
$ws.Range("G8").Value = "KKKK"

# Turn the e-mail address into a real mailto: hyperlink, like the other
# rows above it.
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:kkiwi@example.com", "", "", "kkiwi@example.com")

# Hyperlinks.Add() stamps the cell with Excel's built-in blue/underlined
# "Hyperlink" style; restore the plain style used by the rest of column B.
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Reset the view back to A1 / last-entered cell, matching how the sheet
# was left after the new row was added.
$ws.Range("A1").Select()
$ws.Range("F19").Select()
